# "Generate Report for Archive"
# - Update the "Status" text from "Ready for handoff" to "In Translation"
#   on every sheet that shows it (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# - Narrow the "Status" column(s) (Overview columns E & F, zh-cn column C,
#   de-de column C) to their new, tighter auto-fit width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C holds the status text ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C holds the status text ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

$dede.Columns.Item(3).ColumnWidth = 12.5

Write-Output "Report regenerated for archive"
